$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = '66.259.81'
$ws.Range("E2").Value = '  +7.60%  '

$ws.Range("D3").Value = '3.019.58'
$ws.Range("E3").Value = '  +4.75%  '

$ws.Range("E4").Value = '  +0.14%  '

$ws.Range("D5").Value = '''583.15'
$ws.Range("E5").Value = '  +3.00%  '

$ws.Range("D6").Value = '''157.61'
$ws.Range("E6").Value = '  +10.77%  '

$ws.Range("D7").Value = '''0.999'
$ws.Range("E7").Value = '  -0.08%  '

$ws.Range("D8").Value = '3.015.71'
$ws.Range("E8").Value = '  +4.65%  '

$ws.Range("E10").Value = '  +1.84%  '

$ws.Range("E11").Value = '  +7.53%  '

$ws.Range("E12").Value = '  +5.81%  '

$ws.Range("D13").Value = '''0.0000253'
$ws.Range("E13").Value = '  +10.18%  '

$ws.Range("D14").Value = '''34.66'
$ws.Range("E14").Value = '  +9.80%  '

$ws.Range("E15").Value = '  +0.68%  '

$ws.Range("D16").Value = '66.259.91'
$ws.Range("E16").Value = '  +7.68%  '

$ws.Range("D17").Value = '3.518.74'
$ws.Range("E17").Value = '  +4.71%  '

$ws.Range("D18").Value = '''6.96'
$ws.Range("E18").Value = '  +7.35%  '

$ws.Range("D19").Value = '3.023.21'
$ws.Range("E19").Value = '  +4.70%  '

$ws.Range("D20").Value = '''463.86'
$ws.Range("E20").Value = '  +8.18%  '

$ws.Range("E21").Value = '  +7.16%  '

$ws.Range("D22").Value = '''0.685'
$ws.Range("E22").Value = '  +5.49%  '

$ws.Range("D23").Value = '''7.35'
$ws.Range("E23").Value = '  +8.57%  '

$ws.Range("D24").Value = '''82.32'
$ws.Range("E24").Value = '  +4.50%  '

$ws.Range("E25").Value = '  +12.95%  '

$ws.Range("D26").Value = '''12.49'
$ws.Range("E26").Value = '  +5.98%  '

$ws.Range("D27").Value = '''10.65'
$ws.Range("E27").Value = '  +7.44%  '

$ws.Range("D28").Value = '''1.00'
$ws.Range("E28").Value = '  -0.03%  '

$ws.Range("D29").Value = '''8.05'
$ws.Range("E29").Value = '  +14.69%  '

$ws.Range("D30").Value = '''2.37'
$ws.Range("E30").Value = '  +17.65%  '

$ws.Range("E31").Value = '  +1.16%  '

$ws.Range("E32").Value = '  +5.11%  '

$ws.Range("E33").Value = '  +6.84%  '

$ws.Range("E34").Value = '  +5.38%  '

$ws.Range("D35").Value = '''1.00'
$ws.Range("E35").Value = '  -0.02%  '

$ws.Range("E36").Value = '  +4.30%  '

$ws.Range("D37").Value = '''5.80'
$ws.Range("E37").Value = '  +8.60%  '

$ws.Range("D38").Value = '''2.17'
$ws.Range("E38").Value = '  +14.75%  '

$ws.Range("E39").Value = '  +9.76%  '

$ws.Range("D40").Value = '''49.52'
$ws.Range("E40").Value = '  +1.57%  '

$ws.Range("E41").Value = '  +8.64%  '

$ws.Range("B42").Value = 'TheGraph'
$ws.Range("C42").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D42").Value = '''0.302'
$ws.Range("E42").Value = '  +14.33%  '

$ws.Range("B43").Value = 'Arweave'
$ws.Range("C43").Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range("D43").Value = '''43.81'
$ws.Range("E43").Value = '  +11.30%  '

$ws.Range("D44").Value = '''8.44'
$ws.Range("E44").Value = '  +3.79%  '

$ws.Range("D45").Value = '''390.98'
$ws.Range("E45").Value = '  +14.34%  '

$ws.Range("D46").Value = '2.807.46'

$ws.Range("E47").Value = '  +6.37%  '

$ws.Range("D48").Value = '''133.94'

$ws.Range("D50").Value = '''23.60'
$ws.Range("E50").Value = '  +10.55%  '

$ws.Range("E51").Value = '  +4.58%  '
